$wb = $excel.ActiveWorkbook

# Add violent crime data update for 2023-03-02

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1032
$ws.Range("J3").Value = 1109
$ws.Range("E4").Value = 1981
$ws.Range("J4").Value = 251
$ws.Range("J5").Value = 83
$ws.Range("J6").Value = 1519
$ws.Range("E7").Value = 25984
$ws.Range("J7").Value = 3994

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 30
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J3").Value = 7
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 261
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 38
$ws.Range("J7").Value = 117
$ws.Range("J10").Value = 28
$ws.Range("J18").Value = 59
$ws.Range("J19").Value = 136
$ws.Range("J23").Value = 36
$ws.Range("J27").Value = 21
$ws.Range("J29").Value = 212
$ws.Range("I31").Value = 261
$ws.Range("J31").Value = 28
$ws.Range("J33").Value = 162
$ws.Range("J37").Value = 140
$ws.Range("J42").Value = 171
$ws.Range("J44").Value = 29
$ws.Range("J50").Value = 22
$ws.Range("J51").Value = 53
$ws.Range("J54").Value = 76
$ws.Range("J57").Value = 15
$ws.Range("E63").Value = 327
$ws.Range("J63").Value = 22
$ws.Range("J67").Value = 145
$ws.Range("J72").Value = 17
$ws.Range("J73").Value = 37
$ws.Range("J76").Value = 62
$ws.Range("I78").Value = 345
$ws.Range("J80").Value = 11
$ws.Range("J83").Value = 91
$ws.Range("J84").Value = 40
$ws.Range("J85").Value = 163
$ws.Range("J86").Value = 16
$ws.Range("J91").Value = 55
$ws.Range("J96").Value = 52
$ws.Range("E101").Value = 25984
$ws.Range("J101").Value = 3994

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 36
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 71
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 22
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 76
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 30
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 56
$ws.Range("J7").Value = 163

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 24
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 171

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J2").Value = 7
$ws.Range("J3").Value = 10
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I4").Value = 42
$ws.Range("I7").Value = 345

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J4").Value = 11
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 41
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 117
